$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Update Runmode column (C3:C7) from "Y" to "N", leaving C2 as "Y"
$ws.Range("C3:C7").Value = "N"

# Update selection to C3 active cell, single cell selection
$ws.Range("C3").Select()
